$wb = $excel.ActiveWorkbook

# The "data" sheet holds the panel gene rows; we refresh its time_taken
# column (F) with the new query timestamps and add a sibling "metadata"
# sheet summarising the panel query itself.
$data = $wb.Worksheets.Item("data")

$data.Range("F2").Value = "2021-10-05 14:21:56.742844"
$data.Range("F3").Value = "2021-10-05 14:21:56.742852"
$data.Range("F4").Value = "2021-10-05 14:21:56.742856"
$data.Range("F5").Value = "2021-10-05 14:21:56.742859"
$data.Range("F6").Value = "2021-10-05 14:21:56.742861"
$data.Range("F7").Value = "2021-10-05 14:21:56.742864"
$data.Range("F8").Value = "2021-10-05 14:21:56.742867"
$data.Range("F9").Value = "2021-10-05 14:21:56.742869"
$data.Range("F10").Value = "2021-10-05 14:21:56.742873"

# Add the new "metadata" tab right after "data" (so it becomes sheetId 2,
# appearing second in the tab strip, matching the workbook.xml ordering).
$meta = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $data)
$meta.Name = "metadata"

# Header row (B1:G1) - columns describing the panel query.
$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

# Data row 2 - the single panel record.
$meta.Range("A2").Value = 0
$meta.Range("B2").Value = "Ovarian cancer pertinent cancer susceptibility"
$meta.Range("C2").Value = 117

# data_version "1.5" must stay text (matches the source inlineStr), not
# become the number 1.5: force a text format just long enough to pin the
# type, then drop the number format again so no stray style lingers on
# the cell's visible formatting.
$dataVersion = $meta.Range("D2")
$dataVersion.NumberFormat = "@"
$dataVersion.Value = "1.5"
$dataVersion.ClearFormats()

$meta.Range("E2").Value = "2021-09-29T13:11:49.995302Z"
$meta.Range("F2").Value = "2021-10-05 14:21:56.739213"
$meta.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/117/?format=json"

# Re-use the existing bold/bordered header style (style index 1, already
# present in styles.xml for the "data" sheet's own header row) instead of
# toggling Font/Borders directly, which would mint a brand-new style.
$data.Range("B1").Copy()
$meta.Range("B1:G1").PasteSpecial(-4122) # xlPasteFormats
$meta.Range("A2").PasteSpecial(-4122)    # xlPasteFormats

Write-Host "Added 'metadata' sheet and refreshed data!F2:F10 timestamps"
